# Apply "new version iq and urgency" updates to the Metadata sheet:
#  - Status changes from "draft" to "active"
#  - Date is refreshed to the new publication timestamp
#  - Case Sensitive value is set to "true"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B6").Value = "active"
$ws.Range("B8").Value = "2024-12-16T14:50:05-03:00"

# "true" would normally auto-convert to a native Boolean when assigned via
# .Value (same as typing it into Excel), but the source cell needs to stay a
# plain text/shared-string cell. Force literal text with a leading
# apostrophe, then restore the original (non quote-prefixed) formatting from
# a neighbouring cell that shares the same style.
$ws.Range("B17").Value = "'true"
$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122)
